# Update the "How many observations per UoO" description on the Metadata
# sheet to reflect that the matched persons are non-pregnant MS cases
# (rather than non-MS pregnancies).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")
$ws.Range("B6").Value = "1 + as many as the matched MS non-pregnancies"
